$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: find the Nth (1-based) occurrence of $searchText, re-scanning the
# live/current document content from the start every time, and return a
# Range collapsed exactly over that occurrence.
# ---------------------------------------------------------------------------
function Find-Occurrence($searchText, $occurrence) {
    $count = 0
    $searchStart = 0
    $docEnd = $d.Content.End
    while ($true) {
        $rng = $d.Range($searchStart, $docEnd)
        $found = $rng.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
        if (-not $found) {
            return $null
        }
        $count = $count + 1
        if ($count -eq $occurrence) {
            return $rng
        }
        $searchStart = $rng.End
    }
}

# ---------------------------------------------------------------------------
# Edit 1: first "(Note the last five ..." note -> "(Note: The last five ..."
#   This is the sentence right before the "Input Parameters:" heading of the
#   readDiaSessions section (the one whose closing ")" lives in a separate,
#   already-existing run that we leave untouched).
# ---------------------------------------------------------------------------
$rng = Find-Occurrence "(Note the last five" 1
$rng.Text = "(Note: The last five"

# ---------------------------------------------------------------------------
# Edit 2: fix typo "ouput" -> "output" (readDiaSessions.R track list line)
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.ClearFormatting()
$rng.Find.Execute("a track list ouput from readDiaSessions.R", $true, $false, $false, $false, $false, $true, 1, $false, "a track list output from readDiaSessions.R", 2) | Out-Null

# ---------------------------------------------------------------------------
# Edit 3: outputColWise description - "home directory" -> "current directory"
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.ClearFormatting()
$rng.Find.Execute("similar to a Diatrack output, in the home directory.", $true, $false, $false, $false, $false, $true, 1, $false, "similar to a Diatrack output, in the current directory.", 2) | Out-Null

# ---------------------------------------------------------------------------
# Edit 4: outputRowWise description - "home directory" -> "current directory"
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.ClearFormatting()
$rng.Find.Execute("similar to an ImageJ output, in the home directory.", $true, $false, $false, $false, $false, $true, 1, $false, "similar to an ImageJ output, in the current directory.", 2) | Out-Null

# ---------------------------------------------------------------------------
# Edit 5: linkSkippedFrames tolerance description
#   "to limit how far the next skipped point can deviate from the last point in "
#   -> "to limit how far the next point after the skip can deviate from the last point in "
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.ClearFormatting()
$rng.Find.Execute("to limit how far the next skipped point can deviate from the last point in ", $true, $false, $false, $false, $false, $true, 1, $false, "to limit how far the next point after the skip can deviate from the last point in ", 2) | Out-Null

# ---------------------------------------------------------------------------
# Edit 6: second "(Note the last five ..." note (the one whose closing ")" is
#   already part of the same run) -> "(Note: The last five ...", and the
#   _GoBack bookmark is relocated here, between "T" and "he" of "The".
#   Because edit 1 already rewrote the first occurrence, searching again for
#   the still-unedited text "(Note the last five" and taking the 1st match
#   now correctly resolves to this second paragraph.
# ---------------------------------------------------------------------------
$rng = Find-Occurrence "(Note the last five" 1
$startPos = $rng.Start
$rng.Text = "(Note: The last five"

# Remove the old bookmark (it currently sits right after "distance tolerance level")
$bm = $d.Bookmarks("_GoBack")
if ($bm.Exists) {
    $bm.Delete()
}

# Insert the bookmark between "T" and "he" of the word "The" we just typed,
# i.e. 8 characters ("(Note: T") after the start of the edited range.
$insertPoint = $startPos + 8
$bmRange = $d.Range($insertPoint, $insertPoint)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

# ---------------------------------------------------------------------------
# Edit 7: append " after the frame skip" after " measured in pixels"
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.ClearFormatting()
$rng.Find.Execute("distance tolerance level measured in pixels", $true, $false, $false, $false, $false, $true, 1, $false, "distance tolerance level measured in pixels after the frame skip", 2) | Out-Null
